$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "30.018.88"
Set-TextCell "E2" "  -0.53%  "

Set-TextCell "D3" "1.875.38"
Set-TextCell "E3" "  -1.67%  "

Set-TextCell "D4" "0.9991"
Set-TextCell "E4" "  +0.05%  "

Set-TextCell "D5" "243.14"
Set-TextCell "E5" "  -3.78%  "

Set-TextCell "D6" "0.9991"
Set-TextCell "E6" "  +0.04%  "

Set-TextCell "E7" "  -3.60%  "

Set-TextCell "D8" "0.2926"
Set-TextCell "E8" "  -3.30%  "

Set-TextCell "D9" "0.06594"
Set-TextCell "E9" "  -3.27%  "

Set-TextCell "D10" "1.876.79"
Set-TextCell "E10" "  -1.60%  "

Set-TextCell "D11" "16.63"
Set-TextCell "E11" "  -4.13%  "

Set-TextCell "D12" "0.07177"
Set-TextCell "E12" "  -2.09%  "

Set-TextCell "D13" "0.6672"
Set-TextCell "E13" "  -3.92%  "

Set-TextCell "D14" "86.19"
Set-TextCell "E14" "  -1.13%  "

Set-TextCell "D15" "4.907"
Set-TextCell "E15" "  -0.37%  "

Set-TextCell "D16" "29.976.57"
Set-TextCell "E16" "  -0.65%  "

Set-TextCell "D17" "0.000007791"
Set-TextCell "E17" "  -5.89%  "

Set-TextCell "D18" "0.9990"
Set-TextCell "E18" "  +0.08%  "

Set-TextCell "D19" "12.78"
Set-TextCell "E19" "  -2.39%  "

Set-TextCell "D20" "2.122.51"

Set-TextCell "D21" "0.9950"
Set-TextCell "E21" "  -0.31%  "

Set-TextCell "D22" "4.768"
Set-TextCell "E22" "  -1.29%  "

Set-TextCell "D23" "5.860"
Set-TextCell "E23" "  +1.69%  "

Set-TextCell "D24" "9.087"
Set-TextCell "E24" "  -3.01%  "

Set-TextCell "D25" "151.71"
Set-TextCell "E25" "  +2.39%  "

Set-TextCell "D26" "143.15"
Set-TextCell "E26" "  +6.05%  "

Set-TextCell "D27" "16.92"
Set-TextCell "E27" "  -1.46%  "

Set-TextCell "D28" "1.893"
Set-TextCell "E28" "  -5.71%  "

Set-TextCell "D29" "1.383"
Set-TextCell "E29" "  -0.97%  "

Set-TextCell "D30" "4.195"
Set-TextCell "E30" "  -2.66%  "

Set-TextCell "D31" "0.08738"
Set-TextCell "E31" "  -2.00%  "

Set-TextCell "D32" "3.977"
Set-TextCell "E32" "  -0.88%  "

Set-TextCell "D33" "0.05004"
Set-TextCell "E33" "  -1.72%  "

Set-TextCell "D34" "0.7194"
Set-TextCell "E34" "  -0.38%  "

Set-TextCell "D35" "1.110"
Set-TextCell "E35" "  -3.60%  "

Set-TextCell "D36" "2.660"
Set-TextCell "E36" "  -0.95%  "

Set-TextCell "D37" "0.01821"
Set-TextCell "E37" "  +7.62%  "

Set-TextCell "D38" "2.686"
Set-TextCell "E38" "  -4.77%  "

Set-TextCell "D39" "2.152"
Set-TextCell "E39" "  -6.39%  "

Set-TextCell "D40" "0.9303"
Set-TextCell "E40" "  -3.21%  "

Set-TextCell "B41" "FraxShare"
Set-TextCell "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D41" "5.757"
Set-TextCell "E41" "  -5.48%  "

Set-TextCell "B42" "PaxDollar"
Set-TextCell "C42" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D42" "0.9982"
Set-TextCell "E42" "  -0.09%  "

Set-TextCell "D43" "0.4213"
Set-TextCell "E43" "  -2.72%  "

Set-TextCell "D44" "102.99"
Set-TextCell "E44" "  -2.41%  "

Set-TextCell "D45" "7.357"
Set-TextCell "E45" "  -4.40%  "

Set-TextCell "D46" "0.1268"
Set-TextCell "E46" "  -1.15%  "

Set-TextCell "D47" "0.05698"
Set-TextCell "E47" "  -1.16%  "

Set-TextCell "D48" "32.74"
Set-TextCell "E48" "  -2.46%  "

Set-TextCell "D49" "8.259"
Set-TextCell "E49" "  -1.87%  "

Set-TextCell "D50" "0.3758"
Set-TextCell "E50" "  -1.87%  "

Set-TextCell "E51" "  -2.12%  "
